$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a blank row 1 sitting above the real data: the header
# ("name"/"lat"/"lng") was on row 2, the three data rows were 3-5, and a
# stray formatted-but-empty cell sat at B8. Remove that leading blank row so
# everything shifts up by one: header -> row 1, data -> rows 2-4, stray cell
# -> B7.
$shp = $ws.Shapes("TextBox 1")

# Capture the text box's anchor (row 19, 1-indexed) and its exact vertical
# offset into that row, in EMU, before the sheet is touched. (81643 EMU ~=
# 6.43pt is the offset already baked into the file for this shape.)
$anchorRow = 19
$emuPerPoint = 12700
$offsetIntoRowEmu = 81643

$ws.Rows("1:1").Delete()

# The text box is anchored at a fixed absolute position ("don't move or size
# with cells"), so it does NOT automatically follow the row deletion the way
# cell contents do. Move it back onto the same logical row, now one row
# higher (19 -> 18), using the row's exact top (from cell geometry, which is
# reported at full precision) plus the shape's original offset into the row
# so the anchor offset is reproduced exactly rather than through the
# 1/100-pt-rounded Shape.Top getter.
$newRowTop = $ws.Cells.Item($anchorRow - 1, 1).Top
$shp.Top = $newRowTop + ($offsetIntoRowEmu / $emuPerPoint)

# Leave the selection where the user clicked after making the edit.
$ws.Range("B6").Select() | Out-Null
